$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates ---
# Excel's ColumnWidth setter bakes in a ~0.8333 char padding offset when it
# round-trips to the OOXML <col width=".."> attribute, so subtract that back
# out to land exactly on the target integer widths from the diff.
$pad = 0.8333333333333333
$ws.Columns.Item(2).ColumnWidth = 14 - $pad
$ws.Columns.Item(5).ColumnWidth = 14 - $pad
$ws.Columns.Item(6).ColumnWidth = 11 - $pad
$ws.Columns.Item(7).ColumnWidth = 48 - $pad
$ws.Columns.Item(8).ColumnWidth = 15 - $pad
$ws.Columns.Item(9).ColumnWidth = 30 - $pad
$ws.Columns.Item(10).ColumnWidth = 16 - $pad

# --- Keep rows 5:9 present (unchanged, empty) so they are not dropped ---
foreach ($r in 5..9) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# --- Row 11: rebuild header row with new columns/values, clear styling ---
$ws.Range("A11:J11").ClearFormats()

$ws.Range("A11").Value = "adapter-driver"
$ws.Range("B11").Value = "good sum"
$ws.Range("C11").Value = "critical sum"
$ws.Range("D11").Value = "warning sum"
$ws.Range("E11").Value = "client count"
$ws.Range("F11").Value = "total sum"
$ws.Range("G11").Value = "adapter"
$ws.Range("H11").Value = "driver"
$ws.Range("I11").Value = "good roaming calculation (%)"
$ws.Range("J11").Value = "driver vintage"

# --- Row 12: new data row ---
$ws.Range("A12").Value = "Realtek RTL8852AE WiFi 6 802.11ax PCIe Adapter - 6001.10.356.1"
$ws.Range("B12").Value = 1071383
$ws.Range("C12").Value = 4419
$ws.Range("D12").Value = 180
$ws.Range("E12").Value = 1644
$ws.Range("F12").Value = 1075982
$ws.Range("G12").Value = "realtek rtl8852ae wifi 6 802.11ax pcie adapter"
$ws.Range("H12").Value = "6001.10.356.1"
$ws.Range("I12").Value = 99.59999999999999

# J12 looks like a date ("2024-05-12"); force text so Excel doesn't
# auto-convert it to a date serial, then strip the leftover number-format
# style so the cell keeps the default (unstyled) appearance.
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "2024-05-12"
$ws.Range("J12").ClearFormats()

# --- Remove now-unused rows 13:16 (dimension shrinks to J12) ---
$ws.Range("A13:J16").Delete()
